$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R3").Copy()
$ws.Range("S3").PasteSpecial(-4122)

$ws.Range("S4").Value = 2022
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)

$ws.Range("S5").Value = 135
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)

$ws.Range("S6").Value = 99
$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)

$ws.Range("S7").Value = 36
$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial(-4122)

$ws.Range("S8").Value = 97
$ws.Range("R8").Copy()
$ws.Range("S8").PasteSpecial(-4122)

$ws.Range("S9").Value = 80
$ws.Range("R9").Copy()
$ws.Range("S9").PasteSpecial(-4122)

$ws.Range("S10").Value = 17
$ws.Range("R10").Copy()
$ws.Range("S10").PasteSpecial(-4122)

$ws.Range("S11").Value = 17
$ws.Range("R11").Copy()
$ws.Range("S11").PasteSpecial(-4122)

$ws.Range("S12").Value = 11
$ws.Range("R12").Copy()
$ws.Range("S12").PasteSpecial(-4122)

$ws.Range("S13").Value = 6
$ws.Range("R13").Copy()
$ws.Range("S13").PasteSpecial(-4122)

$ws.Range("S14").Value = 5
$ws.Range("R14").Copy()
$ws.Range("S14").PasteSpecial(-4122)

$ws.Range("S15").Value = 3
$ws.Range("R15").Copy()
$ws.Range("S15").PasteSpecial(-4122)

$ws.Range("S16").Value = 2
$ws.Range("R16").Copy()
$ws.Range("S16").PasteSpecial(-4122)

$ws.Range("S17").Value = "-"
$ws.Range("R17").Copy()
$ws.Range("S17").PasteSpecial(-4122)

$ws.Range("S18").Value = "-"
$ws.Range("R18").Copy()
$ws.Range("S18").PasteSpecial(-4122)

$ws.Range("S19").Value = "-"
$ws.Range("R19").Copy()
$ws.Range("S19").PasteSpecial(-4122)

$ws.Range("S20").Value = 6
$ws.Range("R20").Copy()
$ws.Range("S20").PasteSpecial(-4122)

$ws.Range("S21").Value = 1
$ws.Range("R21").Copy()
$ws.Range("S21").PasteSpecial(-4122)

$ws.Range("S22").Value = 5
$ws.Range("R22").Copy()
$ws.Range("S22").PasteSpecial(-4122)

$ws.Range("S23").Value = "-"
$ws.Range("R23").Copy()
$ws.Range("S23").PasteSpecial(-4122)

$ws.Range("S24").Value = "-"
$ws.Range("R24").Copy()
$ws.Range("S24").PasteSpecial(-4122)

$ws.Range("S25").Value = "-"
$ws.Range("R25").Copy()
$ws.Range("S25").PasteSpecial(-4122)

$ws.Range("S26").Value = 10
$ws.Range("R26").Copy()
$ws.Range("S26").PasteSpecial(-4122)

$ws.Range("S27").Value = 4
$ws.Range("R27").Copy()
$ws.Range("S27").PasteSpecial(-4122)

$ws.Range("S28").Value = 6
$ws.Range("R28").Copy()
$ws.Range("S28").PasteSpecial(-4122)

$ws.Range("S29").Value = "-"
$ws.Range("R29").Copy()
$ws.Range("S29").PasteSpecial(-4122)

$ws.Range("S30").Value = "-"
$ws.Range("R30").Copy()
$ws.Range("S30").PasteSpecial(-4122)

$ws.Range("S31").Value = "-"
$ws.Range("R31").Copy()
$ws.Range("S31").PasteSpecial(-4122)

$ws.Range("S32").Value = "-"
$ws.Range("R32").Copy()
$ws.Range("S32").PasteSpecial(-4122)

$ws.Range("S33").Value = "-"
$ws.Range("R33").Copy()
$ws.Range("S33").PasteSpecial(-4122)

$ws.Range("S34").Value = "-"
$ws.Range("R34").Copy()
$ws.Range("S34").PasteSpecial(-4122)

$ws.Range("T24").Select()
